# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text block on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value()
$text = $text.Replace("1000 Bs = 4.24 = 16441.9 pesos", "1000 Bs = 4.15 = 16058.16 pesos")
$text = $text.Replace("16441.9 pesos = 4.23 = 959.59 Bs", "16058.16 pesos = 4.12 = 930.5 Bs")
$cell.Value = $text

# --- Update the rate table numbers on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 240.999
$wsTasas.Range("O10").Value = 3870
$wsTasas.Range("N12").Value = 3893.3
$wsTasas.Range("O12").Value = 225.6
